$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Extend the formatting that already exists on rows 63/64 down onto the
# newly added rows (65-70), mirroring how these were styled by hand.
$ws.Range("A64:D64").Copy()
$ws.Range("A65:D65").PasteSpecial(-4122)

$ws.Range("A63:D63").Copy()
$ws.Range("A66:D70").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# New problems added to the diary
$data = @(
    @(119, "Pascal's Triangle - II"),
    @(121, "Best Time to Buy and Sell Stock"),
    @(122, "Best Time to Buy and Sell Stock - II"),
    @(125, "Valid Palindrome"),
    @(136, "Single Number"),
    @(141, "Linked List Cycle")
)

$row = 65
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}

# Restore the view state to match what was left after the edit session
$ws.Application.ActiveWindow.ScrollRow = 49
$ws.Range("G64").Select()
